# Daily attendance processing - 2025-11-21 14:47:43
#
# Normalises the "Recorded By" (column G) cell text on the attendance
# report so the listed recorders appear in a consistent order
# (e.g. "System" is listed before a user email, "admin@admin.com"
# before "dnasr281@gmail.com", etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of the exact "Recorded By" text as it previously appeared -> the
# reordered text it should now read.
$map = @{
    "System, system, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
        $changed = $changed + 1
    }
}

Write-Output "Recorded By cells reordered: $changed"
